$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.054.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "'2.759.81"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'578.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'158.53"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'5.77"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -13.94%  "
$ws.Range("D11").Value = "'0.387"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "'3.246.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'27.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").Value = "'63.722.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'2.762.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "'12.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'360.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'65.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "'8.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").Value = "'166.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").Value = "'20.34"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").Value = "'4.94"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "'6.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +13.00%  "
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "'330.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "'21.57"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'21.87"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'0.0595"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").Value = "'0.0258"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "'0.635"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "'136.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +0.60%  "
